$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.350.35"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.872.05"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4671"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2845"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "1.882.71"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.159"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6777"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").Value = "30.336.22"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.523"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("D21").Value = "2.131.71"
$ws.Range("E21").Value = "  -0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007282"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.199"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.321"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.930"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.353"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09692"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.446"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.476"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7041"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.24%  "

$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8497"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4201"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.242"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.303"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "932.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1136"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
